$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) currently lists periods in descending
# order from row 16 (2003) down to row 39 (1804). Update the database so
# that the periods run in ascending order instead (1804 at row 16 up to
# 2003 at row 39), keeping every other cell/style untouched.

$periods = @("1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}
